# "Generate Report for Handoff"
#
# The "b.md" row in every sheet moves from the previous handback status to a
# fresh "Ready for handoff" state, with a new handoff package
# (b.63290e5768f688058c7b37413b0a5c26c308f864.*.xlf) and a new handoff
# timestamp.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: row 3 is the "b.md" file ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = "Ready for handoff"          # zh-cn status
$wsOverview.Range("C3").Value = "Ready for handoff"          # de-de status
$wsOverview.Range("D3").Value = "2016-35-13 08:35:19"        # Latest Handoff Date

# --- zh-cn sheet: row 3 is the "b.md" source file ---
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C3").Value = "Ready for handoff"
$wsZh.Range("D3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$wsZh.Range("E3").Value = "2016-03-13 08:35:15"
foreach ($hl in $wsZh.Hyperlinks) {
    if ($hl.Range.Address() -eq '$D$3') {
        $hl.TextToDisplay = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
    }
}

# --- de-de sheet: row 3 is the "b.md" source file ---
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C3").Value = "Ready for handoff"
$wsDe.Range("D3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$wsDe.Range("E3").Value = "2016-03-13 08:35:19"
foreach ($hl in $wsDe.Hyperlinks) {
    if ($hl.Range.Address() -eq '$D$3') {
        $hl.TextToDisplay = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
    }
}
